# Applies the "LinuxForHealth" rebrand + version/date bump edit described by the diff.
$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aggregate-claim-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Fixed Value of the Extension.url row picks up the same URL rebrand
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aggregate-claim-indicator"

# The ele-1/ext-1 constraint text was incorrectly attached to the "Extension"
# row (AI2); it actually belongs to the "Extension.extension" row (AI4).
$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
